# Apply updated probability values to Sheet1 (team specific time matrix)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1952054794520548
$ws.Range("C2").Value = 0.5582191780821918
$ws.Range("J2").Value = 0.01027397260273973
$ws.Range("P2").Value = 0.1712328767123288
$ws.Range("S2").Value = 0.06506849315068493
$ws.Range("B3").Value = 0.005780346820809248
$ws.Range("C3").Value = 0.03468208092485549
$ws.Range("J3").Value = 0.02890173410404624
$ws.Range("P3").Value = 0.7456647398843931
$ws.Range("S3").Value = 0.1849710982658959
$ws.Range("J4").Value = 0.02040816326530612
$ws.Range("O4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.673469387755102
$ws.Range("S4").Value = 0.2857142857142857
$ws.Range("J5").Value = 0.2
$ws.Range("P5").Value = 0.6
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.07027027027027027
$ws.Range("D6").Value = 0.01081081081081081
$ws.Range("F6").Value = 0.05405405405405406
$ws.Range("J6").Value = 0.2594594594594595
$ws.Range("O6").Value = 0.01081081081081081
$ws.Range("Q6").Value = 0.1837837837837838
$ws.Range("R6").Value = 0.08108108108108109
$ws.Range("S6").Value = 0.3297297297297297
$ws.Range("B7").Value = 0.1219512195121951
$ws.Range("D7").Value = 0.01219512195121951
$ws.Range("E7").Value = 0.006097560975609756
$ws.Range("F7").Value = 0.04878048780487805
$ws.Range("J7").Value = 0.1524390243902439
$ws.Range("O7").Value = 0.02439024390243903
$ws.Range("Q7").Value = 0.1951219512195122
$ws.Range("R7").Value = 0.07926829268292683
$ws.Range("S7").Value = 0.3597560975609756
$ws.Range("B8").Value = 0.09565217391304348
$ws.Range("D8").Value = 0.01391304347826087
$ws.Range("E8").Value = 0.003478260869565218
$ws.Range("F8").Value = 0.04173913043478261
$ws.Range("J8").Value = 0.1234782608695652
$ws.Range("O8").Value = 0.01913043478260869
$ws.Range("Q8").Value = 0.2052173913043478
$ws.Range("R8").Value = 0.09217391304347826
$ws.Range("S8").Value = 0.4052173913043478
$ws.Range("B9").Value = 0.1091703056768559
$ws.Range("D9").Value = 0.02620087336244541
$ws.Range("F9").Value = 0.05676855895196507
$ws.Range("J9").Value = 0.1135371179039301
$ws.Range("O9").Value = 0.01746724890829694
$ws.Range("Q9").Value = 0.2358078602620087
$ws.Range("R9").Value = 0.06986899563318777
$ws.Range("S9").Value = 0.37117903930131
$ws.Range("B10").Value = 0.08207934336525308
$ws.Range("D10").Value = 0.02120383036935704
$ws.Range("E10").Value = 0.002051983584131327
$ws.Range("F10").Value = 0.05129958960328317
$ws.Range("J10").Value = 0.1299589603283174
$ws.Range("O10").Value = 0.01299589603283174
$ws.Range("Q10").Value = 0.2387140902872777
$ws.Range("R10").Value = 0.07797537619699042
$ws.Range("S10").Value = 0.3837209302325582
$ws.Range("G11").Value = 0.1590909090909091
$ws.Range("J11").Value = 0.1060606060606061
$ws.Range("K11").Value = 0.2007575757575757
$ws.Range("L11").Value = 0.5227272727272727
$ws.Range("S11").Value = 0.01136363636363636
$ws.Range("G12").Value = 0.7318840579710145
$ws.Range("J12").Value = 0.1811594202898551
$ws.Range("K12").Value = 0.02173913043478261
$ws.Range("L12").Value = 0.03623188405797102
$ws.Range("S12").Value = 0.02898550724637681
$ws.Range("G13").Value = 0.5227272727272727
$ws.Range("J13").Value = 0.3863636363636364
$ws.Range("S13").Value = 0.09090909090909091
$ws.Range("F15").Value = 0.008403361344537815
$ws.Range("H15").Value = 0.1764705882352941
$ws.Range("I15").Value = 0.04621848739495799
$ws.Range("J15").Value = 0.3949579831932773
$ws.Range("K15").Value = 0.05882352941176471
$ws.Range("M15").Value = 0.01260504201680672
$ws.Range("O15").Value = 0.05042016806722689
$ws.Range("S15").Value = 0.2521008403361344
$ws.Range("F16").Value = 0.009569377990430622
$ws.Range("H16").Value = 0.1913875598086124
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.430622009569378
$ws.Range("K16").Value = 0.09569377990430622
$ws.Range("M16").Value = 0.01913875598086124
$ws.Range("O16").Value = 0.03827751196172249
$ws.Range("S16").Value = 0.1244019138755981
$ws.Range("F17").Value = 0.01549053356282272
$ws.Range("H17").Value = 0.2426850258175559
$ws.Range("I17").Value = 0.08777969018932874
$ws.Range("J17").Value = 0.4010327022375215
$ws.Range("K17").Value = 0.06712564543889846
$ws.Range("M17").Value = 0.02237521514629948
$ws.Range("O17").Value = 0.05507745266781412
$ws.Range("S17").Value = 0.108433734939759
$ws.Range("F18").Value = 0.01428571428571429
$ws.Range("H18").Value = 0.2047619047619048
$ws.Range("I18").Value = 0.06666666666666667
$ws.Range("J18").Value = 0.4428571428571428
$ws.Range("K18").Value = 0.1095238095238095
$ws.Range("M18").Value = 0.0380952380952381
$ws.Range("O18").Value = 0.06190476190476191
$ws.Range("S18").Value = 0.06190476190476191
$ws.Range("F19").Value = 0.01032448377581121
$ws.Range("H19").Value = 0.2323008849557522
$ws.Range("I19").Value = 0.09955752212389381
$ws.Range("J19").Value = 0.3864306784660767
$ws.Range("K19").Value = 0.07964601769911504
$ws.Range("M19").Value = 0.01327433628318584
$ws.Range("N19").Value = 0.002212389380530973
$ws.Range("O19").Value = 0.07669616519174041
$ws.Range("S19").Value = 0.09955752212389381
